# Update "想去人数" (number interested) values in column F, rows 2-9,
# on both the "展览" and "全部类型" worksheets.

$wb = $excel.ActiveWorkbook

$values = @{
    2 = 355
    3 = 98
    4 = 1540
    5 = 22
    6 = 52
    7 = 135
    8 = 57
    9 = 376
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $values.Keys) {
        $ws.Cells.Item($row, 6).Value = $values[$row]
    }
}
